$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3924.375
$ws.Range("I15").Value = 3924.375
$ws.Range("K15").Value = 11773.125
$ws.Range("M15").Value = -11604.125
$ws.Range("H32").Value = 5437
$ws.Range("I32").Value = 3374.5
$ws.Range("J32").Value = 7499.5
$ws.Range("K32").Value = 3374.5
$ws.Range("L32").Value = 7499.5
$ws.Range("M32").Value = -3048.5
$ws.Range("N32").Value = -8151.5
$ws.Range("H106").Value = 22962.666
$ws.Range("I106").Value = 22962.666
$ws.Range("K106").Value = 22962.666
$ws.Range("M106").Value = -22331.666
$ws.Range("H112").Value = 68600.07000000001
$ws.Range("J112").Value = 73403.64
$ws.Range("L112").Value = 220210.92
$ws.Range("N112").Value = -222426.92
$ws.Range("H116").Value = 1181609.1
$ws.Range("I116").Value = 2827422
$ws.Range("J116").Value = 6028.4287
$ws.Range("K116").Value = 2827422
$ws.Range("L116").Value = 6028.4287
$ws.Range("M116").Value = -2823980
$ws.Range("N116").Value = -12912.4287
$ws.Range("H138").Value = 3571.7869
$ws.Range("I138").Value = 615.25
$ws.Range("J138").Value = 4295.837
$ws.Range("K138").Value = 1845.75
$ws.Range("L138").Value = 12887.511
$ws.Range("M138").Value = 3294.25
$ws.Range("N138").Value = -23167.511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2430.8333
$ws.Range("I2").Value = 1892
$ws.Range("K2").Value = 1892
$ws.Range("M2").Value = -1779
$ws.Range("H32").Value = 2494.6853
$ws.Range("I32").Value = 2557.6667
$ws.Range("J32").Value = 1424
$ws.Range("K32").Value = 2557.6667
$ws.Range("L32").Value = 1424
$ws.Range("M32").Value = -2270.6667
$ws.Range("N32").Value = -1998
$ws.Range("H45").Value = 99479.52
$ws.Range("I45").Value = 146274.42
$ws.Range("K45").Value = 146274.42
$ws.Range("M45").Value = -145897.42
$ws.Range("H46").Value = 4721.9
$ws.Range("I46").Value = 999.6667
$ws.Range("J46").Value = 6317.143
$ws.Range("K46").Value = 999.6667
$ws.Range("L46").Value = 6317.143
$ws.Range("M46").Value = -680.6667
$ws.Range("N46").Value = -6955.143
$ws.Range("H116").Value = 2430.8333
$ws.Range("I116").Value = 1892
$ws.Range("K116").Value = 1892
$ws.Range("M116").Value = 402
$ws.Range("H132").Value = 1769.6097
$ws.Range("I132").Value = 997.2727
$ws.Range("J132").Value = 4955.5
$ws.Range("K132").Value = 2991.8181
$ws.Range("L132").Value = 14866.5
$ws.Range("M132").Value = -461.8181
$ws.Range("N132").Value = -19926.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2430.8333
$ws.Range("I3").Value = 1892
$ws.Range("K3").Value = 1892
$ws.Range("M3").Value = -1778
$ws.Range("H105").Value = 75112.14
$ws.Range("I105").Value = 101757.2
$ws.Range("K105").Value = 101757.2
$ws.Range("M105").Value = -100010.2
$ws.Range("H107").Value = 3808.7778
$ws.Range("I107").Value = 3830
$ws.Range("K107").Value = 3830
$ws.Range("M107").Value = -1910
$ws.Range("H134").Value = 2384.348
$ws.Range("I134").Value = 1540.5641
$ws.Range("J134").Value = 7085.4287
$ws.Range("K134").Value = 4621.692300000001
$ws.Range("L134").Value = 21256.2861
$ws.Range("M134").Value = -2086.692300000001
$ws.Range("N134").Value = -26326.2861
$ws.Range("H135").Value = 93848.336
$ws.Range("J135").Value = 93848.336
$ws.Range("L135").Value = 93848.336
$ws.Range("N135").Value = -103988.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 29166.666
$ws.Range("I36").Value = 18750
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 18750
$ws.Range("L36").Value = 50000
$ws.Range("M36").Value = -18362
$ws.Range("N36").Value = -50776
$ws.Range("H40").Value = 29166.666
$ws.Range("I40").Value = 18750
$ws.Range("J40").Value = 50000
$ws.Range("K40").Value = 18750
$ws.Range("L40").Value = 50000
$ws.Range("M40").Value = -18590
$ws.Range("N40").Value = -50320
$ws.Range("H54").Value = 32500
$ws.Range("J54").Value = 32500
$ws.Range("L54").Value = 32500
$ws.Range("N54").Value = -33816
$ws.Range("H99").Value = 7260182
$ws.Range("I99").Value = 11613962
$ws.Range("J99").Value = 3882.3333
$ws.Range("K99").Value = 11613962
$ws.Range("L99").Value = 3882.3333
$ws.Range("M99").Value = -11612464
$ws.Range("N99").Value = -6878.3333
$ws.Range("H122").Value = 21142.715
$ws.Range("I122").Value = 26999.8
$ws.Range("K122").Value = 80999.39999999999
$ws.Range("M122").Value = -78549.39999999999
$ws.Range("H126").Value = 7260182
$ws.Range("I126").Value = 11613962
$ws.Range("J126").Value = 3882.3333
$ws.Range("K126").Value = 34841886
$ws.Range("L126").Value = 11646.9999
$ws.Range("M126").Value = -34839416
$ws.Range("N126").Value = -16586.9999
$ws.Range("H132").Value = 15887.929
$ws.Range("J132").Value = 103277.5
$ws.Range("L132").Value = 309832.5
$ws.Range("N132").Value = -314892.5
$ws.Range("H134").Value = 3417.9778
$ws.Range("I134").Value = 1796.8572
$ws.Range("K134").Value = 5390.571599999999
$ws.Range("M134").Value = -2855.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 150250
$ws.Range("I11").Value = 200000
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 600000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -599860
$ws.Range("N11").Value = -3280
$ws.Range("H69").Value = 3388.4443
$ws.Range("J69").Value = 3785.1428
$ws.Range("L69").Value = 11355.4284
$ws.Range("N69").Value = -12977.4284
$ws.Range("H72").Value = 3388.4443
$ws.Range("J72").Value = 3785.1428
$ws.Range("L72").Value = 34066.2852
$ws.Range("N72").Value = -42178.2852
$ws.Range("H110").Value = 49899
$ws.Range("I110").Value = 49899
$ws.Range("K110").Value = 149697
$ws.Range("M110").Value = -145607
$ws.Range("H130").Value = 16250

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3784.875
$ws.Range("I16").Value = 4932.4165
$ws.Range("K16").Value = 4932.4165
$ws.Range("M16").Value = -4762.4165
$ws.Range("H22").Value = 3048.3572
$ws.Range("I22").Value = 3299.111
$ws.Range("J22").Value = 2597
$ws.Range("K22").Value = 3299.111
$ws.Range("L22").Value = 2597
$ws.Range("M22").Value = -3004.111
$ws.Range("N22").Value = -3187
$ws.Range("H27").Value = 3048.3572
$ws.Range("I27").Value = 3299.111
$ws.Range("J27").Value = 2597
$ws.Range("K27").Value = 3299.111
$ws.Range("L27").Value = 2597
$ws.Range("M27").Value = -3192.111
$ws.Range("N27").Value = -2811
$ws.Range("H32").Value = 1287
$ws.Range("I32").Value = 1287
$ws.Range("K32").Value = 1287
$ws.Range("M32").Value = -970
$ws.Range("H46").Value = 1417.6957
$ws.Range("J46").Value = 1954.5
$ws.Range("L46").Value = 1954.5
$ws.Range("N46").Value = -2330.5
$ws.Range("H61").Value = 14163.728
$ws.Range("I61").Value = 15300.167
$ws.Range("J61").Value = 12800
$ws.Range("K61").Value = 15300.167
$ws.Range("L61").Value = 12800
$ws.Range("M61").Value = -15098.167
$ws.Range("N61").Value = -13204
$ws.Range("H82").Value = 2924.4285
$ws.Range("J82").Value = 2599.5
$ws.Range("L82").Value = 2599.5
$ws.Range("N82").Value = -3321.5
$ws.Range("H85").Value = 2924.4285
$ws.Range("J85").Value = 2599.5
$ws.Range("L85").Value = 2599.5
$ws.Range("N85").Value = -5095.5
$ws.Range("H113").Value = 14163.728
$ws.Range("I113").Value = 15300.167
$ws.Range("J113").Value = 12800
$ws.Range("K113").Value = 15300.167
$ws.Range("L113").Value = 12800
$ws.Range("M113").Value = -13130.167
$ws.Range("N113").Value = -17140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2055.1667
$ws.Range("J113").Value = 5997.5
$ws.Range("L113").Value = 17992.5
$ws.Range("N113").Value = -22332.5
